$wb = $excel.ActiveWorkbook

# --- Sheet "계획표" (plan table) ---
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row above the old row 10 ("탐색"), pushing 탐색/그리디/DP down one row.
$ws1.Rows.Item(10).Insert()

# The row-insert duplicates row 9's formatting across columns C:H for the new
# row 10; only A10/B10 should carry content/format, so clear the rest back out.
$ws1.Range("C10:H10").Clear()

# New data point added alongside the existing 240522 entry in row 9.
$ws1.Range("D9").Value = 240523
$ws1.Range("D9").HorizontalAlignment = -4108
$ws1.Range("D9").VerticalAlignment = -4108

# New BFS study entry in the freshly inserted row.
$ws1.Range("B10").Value = "2178 (BFS)"

# The "DFS/BFS" label now spans both rows.
$ws1.Range("A9:A10").Merge()
$ws1.Range("A9:A10").HorizontalAlignment = -4108
$ws1.Range("A9:A10").VerticalAlignment = -4108

# --- Sheet "빈출유형" (frequent types) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B5").Value = "스택/큐"

# Restore selections (사heet2 touched last moves the active tab there, so
# reselect sheet1 afterwards to keep it the active tab).
$ws2.Activate()
$ws2.Range("B13").Select()

$ws1.Activate()
$ws1.Range("C10").Select()
